$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.2332
$ws.Range("C3").Value = -12.2611
$ws.Range("E3").Value = 15.70700000000001
$ws.Range("E12").Value = 17.41800000000001
$ws.Range("C14").Value = -12.53
$ws.Range("C16").Value = -13.3947
$ws.Range("D18").Value = -8.994400000000002
$ws.Range("C21").Value = -12.1927
$ws.Range("C23").Value = -12.51320000000001
$ws.Range("D24").Value = -7.210000000000008
$ws.Range("E24").Value = 16.89910000000001
$ws.Range("C25").Value = -12.569
$ws.Range("D25").Value = -7.839000000000006
$ws.Range("E25").Value = 16.77440000000001
$ws.Range("C26").Value = -12.91970000000001
$ws.Range("D27").Value = -8.923700000000007
$ws.Range("C29").Value = -10.57700000000002
$ws.Range("D30").Value = -7.594200000000005
$ws.Range("D31").Value = -8.615300000000008
$ws.Range("D39").Value = -8.0891
$ws.Range("C40").Value = -13.6774
$ws.Range("E41").Value = 16.12649999999999
$ws.Range("D42").Value = -8.559199999999995
$ws.Range("D48").Value = -7.401600000000001
$ws.Range("E50").Value = 16.4745
$ws.Range("D51").Value = -7.919999999999998
$ws.Range("D52").Value = -7.4873
$ws.Range("C53").Value = -10.67240000000001
$ws.Range("E53").Value = 16.94570000000001
$ws.Range("D55").Value = -8.941500000000003
$ws.Range("D56").Value = -8.139999999999995
$ws.Range("E56").Value = 16.71360000000001
$ws.Range("C57").Value = -14.0552
$ws.Range("D57").Value = -8.190799999999994
$ws.Range("E57").Value = 16.62880000000001
$ws.Range("E58").Value = 16.16110000000002
$ws.Range("C59").Value = -12.9263
$ws.Range("D60").Value = -8.200799999999997
$ws.Range("E61").Value = 16.60520000000001
$ws.Range("E63").Value = 17.51410000000001
$ws.Range("E64").Value = 17.45050000000001
$ws.Range("C65").Value = -11.88830000000001
$ws.Range("C69").Value = -10.5531
$ws.Range("E70").Value = 17.42210000000001
$ws.Range("E72").Value = 16.87960000000001
$ws.Range("D73").Value = -7.499799999999999
$ws.Range("D74").Value = -8.398100000000007
$ws.Range("C79").Value = -10.61410000000001
$ws.Range("C83").Value = -14.28519999999999
$ws.Range("E86").Value = 16.6408
$ws.Range("D89").Value = -7.115399999999993
$ws.Range("E89").Value = 17.43930000000001
$ws.Range("D90").Value = -7.895500000000006
$ws.Range("C91").Value = -10.0292
$ws.Range("D92").Value = -5.900400000000002
$ws.Range("C93").Value = -11.47400000000001
$ws.Range("E98").Value = 15.1905
$ws.Range("C100").Value = -13.0565
$ws.Range("E100").Value = 16.93550000000002
$ws.Range("E102").Value = 16.01049999999999
